$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values B8:AK8 from 20 to 2000
$ws.Range("B8:AK8").Value = 2000

# Update selection / view - scroll back to A1 top-left, select C8:AK8
$ws.Range("C8:AK8").Select()
